$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9
$values = @("F0008", "H", "8", "8", "64", "1", "F", "8", "2", "Aaron", "26/08/23", "20:59")

for ($i = 0; $i -lt $values.Length; $i++) {
    $cell = $ws.Cells.Item($row, $i + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$i]
    $cell.ClearFormats()
}
